$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.077679419419598
$ws.Range("D2").Value = 1.077051058944554
$ws.Range("E2").Value = 1.080426458009901
$ws.Range("F2").Value = 1.089170678772191
$ws.Range("I2").Value = 1.050287125591677
$ws.Range("J2").Value = 1.082572676276498
$ws.Range("K2").Value = 1.079733377003638
$ws.Range("L2").Value = 1.083099922956954
$ws.Range("M2").Value = 1.091821484445301
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.079274077384808
$ws.Range("D3").Value = 1.078302463502725
$ws.Range("E3").Value = 1.081814909288513
$ws.Range("F3").Value = 1.090561679395765
$ws.Range("I3").Value = 1.050673018726119
$ws.Range("J3").Value = 1.083824640796738
$ws.Range("K3").Value = 1.080801342285465
$ws.Range("L3").Value = 1.084305246673493
$ws.Range("M3").Value = 1.093031003619709
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.080304654526588
$ws.Range("D4").Value = 1.079110843163895
$ws.Range("E4").Value = 1.08271230434714
$ws.Range("F4").Value = 1.091460663694051
$ws.Range("I4").Value = 1.050920746271911
$ws.Range("J4").Value = 1.084633051092644
$ws.Range("K4").Value = 1.081490460681232
$ws.Range("L4").Value = 1.085083608082098
$ws.Range("M4").Value = 1.093812017004684
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.080737612399023
$ws.Range("D5").Value = 1.079450365168009
$ws.Range("E5").Value = 1.083089329772963
$ws.Range("F5").Value = 1.091838342704077
$ws.Range("I5").Value = 1.05102442116833
$ws.Range("J5").Value = 1.084972507179559
$ws.Range("K5").Value = 1.081779709941691
$ws.Range("L5").Value = 1.085410462539439
$ws.Range("M5").Value = 1.094139971259655
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.080810290761983
$ws.Range("D6").Value = 1.079507353751883
$ws.Range("E6").Value = 1.083152620136099
$ws.Range("F6").Value = 1.091901741952301
$ws.Range("I6").Value = 1.051041801153961
$ws.Range("J6").Value = 1.085029480092797
$ws.Range("K6").Value = 1.081828249566462
$ws.Range("L6").Value = 1.08546532138539
$ws.Range("M6").Value = 1.094195013879438
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.080310440886965
$ws.Range("D7").Value = 1.07911538112735
$ws.Range("E7").Value = 1.082717343112673
$ws.Range("F7").Value = 1.091465711250012
$ws.Range("I7").Value = 1.050922133424182
$ws.Range("J7").Value = 1.084637588484169
$ws.Range("K7").Value = 1.081494327429915
$ws.Range("L7").Value = 1.08508797696902
$ws.Range("M7").Value = 1.093816400646545
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.078218610023317
$ws.Range("D8").Value = 1.077474262429914
$ws.Range("E8").Value = 1.080895907699891
$ws.Range("F8").Value = 1.089641002472956
$ws.Range("I8").Value = 1.050417949399115
$ws.Range("J8").Value = 1.082996138053326
$ws.Range("K8").Value = 1.080094702859811
$ws.Range("L8").Value = 1.083507594925812
$ws.Range("M8").Value = 1.092230587244521
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.074522424224321
$ws.Range("D9").Value = 1.074571719545457
$ws.Range("E9").Value = 1.077678174874105
$ws.Range("F9").Value = 1.086417042954855
$ws.Range("I9").Value = 1.049514320603024
$ws.Range("J9").Value = 1.080090442810736
$ws.Range("K9").Value = 1.077613400672974
$ws.Range("L9").Value = 1.080710529878792
$ws.Range("M9").Value = 1.089423473702824
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072050994147459
$ws.Range("D10").Value = 1.072629164192416
$ws.Range("E10").Value = 1.075527166485464
$ws.Range("F10").Value = 1.084261580086385
$ws.Range("I10").Value = 1.048901551965604
$ws.Range("J10").Value = 1.078144021445737
$ws.Range("K10").Value = 1.075948800878151
$ws.Range("L10").Value = 1.078837239361449
$ws.Range("M10").Value = 1.087543166998733
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.070978991348958
$ws.Range("D11").Value = 1.071786153617813
$ws.Range("E11").Value = 1.074594281183618
$ws.Range("F11").Value = 1.083326695066478
$ws.Range("I11").Value = 1.048633731171729
$ws.Range("J11").Value = 1.077298909981146
$ws.Range("K11").Value = 1.075225471857245
$ws.Range("L11").Value = 1.078023968494262
$ws.Range("M11").Value = 1.08672678366843
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.070580512906883
$ws.Range("D12").Value = 1.071472734056825
$ws.Range("E12").Value = 1.074247535475342
$ws.Range("F12").Value = 1.08297919601014
$ws.Range("I12").Value = 1.048533874133649
$ws.Range("J12").Value = 1.076984645206498
$ws.Range("K12").Value = 1.074956406354372
$ws.Range("L12").Value = 1.077721557425901
$ws.Range("M12").Value = 1.086423205383225
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.070666001116159
$ws.Range("D13").Value = 1.071539976753687
$ws.Range("E13").Value = 1.074321924162376
$ws.Range("F13").Value = 1.083053746769508
$ws.Range("I13").Value = 1.048555310886746
$ws.Range("J13").Value = 1.077052072150658
$ws.Range("K13").Value = 1.075014139518917
$ws.Range("L13").Value = 1.077786440494266
$ws.Range("M13").Value = 1.086488339325177
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.070946058946694
$ws.Range("D14").Value = 1.071760252171462
$ws.Range("E14").Value = 1.074565623822622
$ws.Range("F14").Value = 1.083297975648484
$ws.Range("I14").Value = 1.048625484653317
$ws.Range("J14").Value = 1.077272939996232
$ws.Range("K14").Value = 1.075203238789112
$ws.Range("L14").Value = 1.077998977793134
$ws.Range("M14").Value = 1.086701696711021
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.071118573245301
$ws.Range("D15").Value = 1.071895932773902
$ws.Range("E15").Value = 1.074715744446885
$ws.Range("F15").Value = 1.083448420944709
$ws.Range("I15").Value = 1.048668671060174
$ws.Range("J15").Value = 1.077408976998336
$ws.Range("K15").Value = 1.075319697367181
$ws.Range("L15").Value = 1.078129885621675
$ws.Range("M15").Value = 1.08683310835218
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.07212209947403
$ws.Range("D16").Value = 1.072685072026312
$ws.Range("E16").Value = 1.075589047091703
$ws.Range("F16").Value = 1.084323591921081
$ws.Range("I16").Value = 1.048919273702947
$ws.Range("J16").Value = 1.078200059653883
$ws.Range("K16").Value = 1.075996751654728
$ws.Range("L16").Value = 1.078891168130565
$ws.Range("M16").Value = 1.08759730077002
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.07275108079777
$ws.Range("D17").Value = 1.073179572478613
$ws.Range("E17").Value = 1.076136444023509
$ws.Range("F17").Value = 1.08487214188839
$ws.Range("I17").Value = 1.049075802139898
$ws.Range("J17").Value = 1.078695664359239
$ws.Range("K17").Value = 1.076420763736315
$ws.Range("L17").Value = 1.079368127036672
$ws.Range("M17").Value = 1.088076065113103
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.073117776431833
$ws.Range("D18").Value = 1.073467826133975
$ws.Range("E18").Value = 1.076455588675531
$ws.Range("F18").Value = 1.085191952239888
$ws.Range("I18").Value = 1.049166862662682
$ws.Range("J18").Value = 1.078984520967246
$ws.Range("K18").Value = 1.076667837557832
$ws.Range("L18").Value = 1.07964612466548
$ws.Range("M18").Value = 1.088355108638813
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073242780181228
$ws.Range("D19").Value = 1.073566082837924
$ws.Range("E19").Value = 1.076564384774239
$ws.Range("F19").Value = 1.085300974160932
$ws.Range("I19").Value = 1.049197871343088
$ws.Range("J19").Value = 1.079082976306127
$ws.Range("K19").Value = 1.076752041937877
$ws.Range("L19").Value = 1.079740880208169
$ws.Range("M19").Value = 1.08845021958468
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.072683615591904
$ws.Range("D20").Value = 1.073126535936783
$ws.Range("E20").Value = 1.076077728329847
$ws.Range("F20").Value = 1.084813303192476
$ws.Range("I20").Value = 1.049059032949969
$ws.Range("J20").Value = 1.078642513601925
$ws.Range("K20").Value = 1.076375296666317
$ws.Range("L20").Value = 1.079316975055973
$ws.Range("M20").Value = 1.088024720168671
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.070863596949502
$ws.Range("D21").Value = 1.071695394553101
$ws.Range("E21").Value = 1.074493866772486
$ws.Range("F21").Value = 1.083226063047526
$ws.Range("I21").Value = 1.048604830642212
$ws.Range("J21").Value = 1.077207909718576
$ws.Range("K21").Value = 1.075147564553886
$ws.Range("L21").Value = 1.077936399913182
$ws.Range("M21").Value = 1.086638877642013
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.069717602744563
$ws.Range("D22").Value = 1.070793910084483
$ws.Range("E22").Value = 1.073496692359458
$ws.Range("F22").Value = 1.08222670330747
$ws.Range("I22").Value = 1.048317076135127
$ws.Range("J22").Value = 1.076303872682771
$ws.Range("K22").Value = 1.074373387423526
$ws.Range("L22").Value = 1.077066487555472
$ws.Range("M22").Value = 1.085765589229425
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.070325278146967
$ws.Range("D23").Value = 1.071271964537073
$ws.Range("E23").Value = 1.074025442723107
$ws.Range("F23").Value = 1.08275661791412
$ws.Range("I23").Value = 1.048469827736632
$ws.Range("J23").Value = 1.076783316137577
$ws.Range("K23").Value = 1.07478400907725
$ws.Range("L23").Value = 1.077527826051141
$ws.Range("M23").Value = 1.086228723501843
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.072714100766299
$ws.Range("D24").Value = 1.073150501423842
$ws.Range("E24").Value = 1.07610425986568
$ws.Range("F24").Value = 1.084839890328897
$ws.Range("I24").Value = 1.049066610966352
$ws.Range("J24").Value = 1.078666530825124
$ws.Range("K24").Value = 1.07639584203617
$ws.Range("L24").Value = 1.079340089066148
$ws.Range("M24").Value = 1.08804792139306
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.075479227241966
$ws.Range("D25").Value = 1.075323397667219
$ws.Range("E25").Value = 1.078511039105214
$ws.Range("F25").Value = 1.087251570768124
$ws.Range("I25").Value = 1.049749744496202
$ws.Range("J25").Value = 1.080843244094678
$ws.Range("K25").Value = 1.0782566856917
$ws.Range("L25").Value = 1.081435124226183
$ws.Range("M25").Value = 1.090150722027585
